$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new weekly record at row 35
$ws.Rows(35).Insert()
$ws.Cells.Item(35,1).Value = 5
$ws.Cells.Item(35,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(35,3).Value = "Maule"
$ws.Cells.Item(35,4).Value = 44818
$ws.Cells.Item(35,5).Value = 7
$ws.Cells.Item(35,6).Value = 300000000
$ws.Cells.Item(35,7).Value = "Espárragos"
$ws.Cells.Item(35,8).Value = "Sin especificar"
$ws.Cells.Item(35,9).Value = "Primera"
$ws.Cells.Item(35,10).Value = 2000
$ws.Cells.Item(35,11).Value = 2800
$ws.Cells.Item(35,12).Value = 2800
$ws.Cells.Item(35,13).Value = 2800
$ws.Cells.Item(35,14).Value = "`$/kilo"
$ws.Cells.Item(35,15).Value = "Provincia de Linares"
$ws.Cells.Item(35,16).Value = 2800
$ws.Cells.Item(35,17).Value = 1
$ws.Cells.Item(35,18).Value = "Hortaliza"

# Insert new weekly record at row 39
$ws.Rows(39).Insert()
$ws.Cells.Item(39,1).Value = 5
$ws.Cells.Item(39,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(39,3).Value = "Maule"
$ws.Cells.Item(39,4).Value = 44824
$ws.Cells.Item(39,5).Value = 7
$ws.Cells.Item(39,6).Value = 300000000
$ws.Cells.Item(39,7).Value = "Espárragos"
$ws.Cells.Item(39,8).Value = "Sin especificar"
$ws.Cells.Item(39,9).Value = "Primera"
$ws.Cells.Item(39,10).Value = 1000
$ws.Cells.Item(39,11).Value = 2800
$ws.Cells.Item(39,12).Value = 2800
$ws.Cells.Item(39,13).Value = 2800
$ws.Cells.Item(39,14).Value = "`$/kilo"
$ws.Cells.Item(39,15).Value = "Provincia de Linares"
$ws.Cells.Item(39,16).Value = 2800
$ws.Cells.Item(39,17).Value = 1
$ws.Cells.Item(39,18).Value = "Hortaliza"

# Insert new weekly record at row 52
$ws.Rows(52).Insert()
$ws.Cells.Item(52,1).Value = 5
$ws.Cells.Item(52,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(52,3).Value = "Maule"
$ws.Cells.Item(52,4).Value = 44816
$ws.Cells.Item(52,5).Value = 7
$ws.Cells.Item(52,6).Value = 300000000
$ws.Cells.Item(52,7).Value = "Espárragos"
$ws.Cells.Item(52,8).Value = "Sin especificar"
$ws.Cells.Item(52,9).Value = "Primera"
$ws.Cells.Item(52,10).Value = 500
$ws.Cells.Item(52,11).Value = 2800
$ws.Cells.Item(52,12).Value = 2800
$ws.Cells.Item(52,13).Value = 2800
$ws.Cells.Item(52,14).Value = "`$/kilo"
$ws.Cells.Item(52,15).Value = "Provincia de Linares"
$ws.Cells.Item(52,16).Value = 2800
$ws.Cells.Item(52,17).Value = 1
$ws.Cells.Item(52,18).Value = "Hortaliza"
